$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# 1. Delete the "Verificat consecventa" / "Scores are total or per round" row.
$t.Rows.Item(4).Delete()

Write-Output "Step 1 done - row deleted"

# 2. Apply green highlight to 4 cells (both pPr/rPr and run rPr pick up
#    <w:highlight w:val="green"/> when we highlight the whole cell range).
$targets = @(
    "Training should also be done with teams reversed",
    "Just for predict by game",
    "Fara numpy si statistics in teamclass",
    "Tournament prediction version"
)
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $row = $t.Rows.Item($i)
    for ($j = 1; $j -le $row.Cells.Count; $j++) {
        $cell = $row.Cells.Item($j)
        $txt = $cell.Range.Text
        foreach ($target in $targets) {
            if ($txt -like "$target*") {
                $cell.Range.Font.HighlightColorIndex = 4
            }
        }
    }
}

Write-Output "Step 2 done - highlights applied"

# 3. Insert three new rows after "Visual interface ?" (and before the
#    trailing blank row), describing the new scoring parameters.
#    Inserting each new row right before the (still) trailing blank row
#    would reverse their order, so insert them in reverse.
$lastRow = $t.Rows.Item($t.Rows.Count)

$row3 = $t.Rows.Add($lastRow)
$row3.Cells.Item(1).Range.Text = "Adugat afisare scor real / scor prezis / diferenta"

$row2 = $t.Rows.Add($lastRow)
$row2.Cells.Item(1).Range.Text = "Adaugat accuracy diferentiat in caz ca nu avem extra si penaltiuri"

$row1 = $t.Rows.Add($lastRow)
$row1.Cells.Item(1).Range.Text = "Adaugat mai multe tipuri de scor`rR1, R2, 90, is_extra, extra, is_penalties, penalties"

Write-Output "Step 3 done - rows inserted"
